# Updates cryptos list price/volume figures (and a few reordered rows) per the
# Sun Oct  8 11:08:16 UTC 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.710.96"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "1.609.05"
$ws.Range("E3").Value = "  -1.94%  "

# Row 4
$ws.Range("D4").Value = "'0.992"
$ws.Range("E4").Value = "  -0.81%  "

# Row 5
$ws.Range("D5").Value = "'208.00"
$ws.Range("E5").Value = "  -2.33%  "

# Row 6
$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  -1.05%  "

# Row 7
$ws.Range("D7").Value = "'0.991"
$ws.Range("E7").Value = "  -0.91%  "

# Row 8
$ws.Range("D8").Value = "'23.09"
$ws.Range("E8").Value = "  -2.06%  "

# Row 9
$ws.Range("E9").Value = "  -2.40%  "

# Row 10
$ws.Range("D10").Value = "'0.0605"
$ws.Range("E10").Value = "  -1.58%  "

# Row 11
$ws.Range("D11").Value = "'0.0872"
$ws.Range("E11").Value = "  -1.07%  "

# Row 12
$ws.Range("D12").Value = "1.831.86"
$ws.Range("E12").Value = "  -2.21%  "

# Row 13
$ws.Range("D13").Value = "1.602.66"
$ws.Range("E13").Value = "  -2.30%  "

# Row 14
$ws.Range("E14").Value = "  -2.92%  "

# Row 15
$ws.Range("D15").Value = "'0.555"
$ws.Range("E15").Value = "  -3.36%  "

# Row 16
$ws.Range("D16").Value = "'64.64"
$ws.Range("E16").Value = "  -1.94%  "

# Row 17
$ws.Range("D17").Value = "27.640.76"
$ws.Range("E17").Value = "  -1.15%  "

# Row 18
$ws.Range("D18").Value = "'227.22"
$ws.Range("E18").Value = "  -2.87%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0714"
$ws.Range("E19").Value = "  -1.35%  "

# Row 20
$ws.Range("D20").Value = "'7.53"
$ws.Range("E20").Value = "  -1.22%  "

# Row 21
$ws.Range("D21").Value = "'0.990"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("D22").Value = "'4.26"
$ws.Range("E22").Value = "  -2.54%  "

# Row 23
$ws.Range("D23").Value = "'10.01"
$ws.Range("E23").Value = "  -6.58%  "

# Row 24
$ws.Range("D24").Value = "'2.01"
$ws.Range("E24").Value = "  -3.32%  "

# Row 25
$ws.Range("D25").Value = "'153.53"
$ws.Range("E25").Value = "  +1.65%  "

# Row 26
$ws.Range("D26").Value = "'6.83"
$ws.Range("E26").Value = "  -2.11%  "

# Row 28
$ws.Range("D28").Value = "'15.34"
$ws.Range("E28").Value = "  -2.16%  "

# Row 29
$ws.Range("D29").Value = "'0.991"
$ws.Range("E29").Value = "  -0.88%  "

# Row 30
$ws.Range("E30").Value = "  -2.52%  "

# Row 31
$ws.Range("E31").Value = "  -1.40%  "

# Row 32
$ws.Range("D32").Value = "'3.37"
$ws.Range("E32").Value = "  +0.52%  "

# Row 33
$ws.Range("D33").Value = "'3.05"
$ws.Range("E33").Value = "  -2.45%  "

# Row 34
$ws.Range("D34").Value = "1.376.40"
$ws.Range("E34").Value = "  -3.39%  "

# Row 36
$ws.Range("D36").Value = "'0.982"
$ws.Range("E36").Value = "  +8.30%  "

# Row 37
$ws.Range("E37").Value = "  -1.93%  "

# Row 38
$ws.Range("D38").Value = "'0.0168"
$ws.Range("E38").Value = "  -0.44%  "

# Row 39
$ws.Range("D39").Value = "'0.550"
$ws.Range("E39").Value = "  -1.38%  "

# Row 40
$ws.Range("D40").Value = "'0.844"
$ws.Range("E40").Value = "  -4.34%  "

# Row 41
$ws.Range("D41").Value = "'1.01"
$ws.Range("E41").Value = "  -1.30%  "

# Row 42
$ws.Range("D42").Value = "'0.990"
$ws.Range("E42").Value = "  -1.02%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "'1.81"
$ws.Range("E43").Value = "  -3.29%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").Value = "'65.06"
$ws.Range("E44").Value = "  -2.06%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.40"
$ws.Range("E45").Value = "  -2.12%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.745.75"
$ws.Range("E46").Value = "  -2.01%  "

# Row 47
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.14"
$ws.Range("E47").Value = "  -2.98%  "

# Row 48
$ws.Range("D48").Value = "'87.09"
$ws.Range("E48").Value = "  -0.88%  "

# Row 49
$ws.Range("E49").Value = "  +0.05%  "

# Row 50
$ws.Range("D50").Value = "'0.0501"
$ws.Range("E50").Value = "  -1.11%  "

# Row 51
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₇0959"
$ws.Range("E51").Value = "  -8.99%  "

